$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows (Índice, Distancia, max, min, Tempo)
$data = @(
    @(0, 3850.2,               4202, 3414, 0.08752556641896565),
    @(1, 3459.533333333333,    3786, 2983, 0.08764753341674805),
    @(2, 3923.9,               4223, 3592, 0.09182977676391602),
    @(3, 3642.4,               3939, 3327, 0.08825164635976156),
    @(4, 3577.066666666667,    3903, 3147, 0.08847372531890869),
    @(5, 3757.633333333333,    4088, 3293, 0.09275384744008382),
    @(6, 4125.133333333333,    4397, 3694, 0.09229811827341715),
    @(7, 3563,                 3869, 3356, 0.0887044350306193),
    @(8, 3881.866666666667,    4150, 3362, 0.08835875193277995),
    @(9, 3789.766666666667,    4217, 3319, 0.08782593409220378)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $rowIndex++
}
